$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block of rows 342:459 (Fecha/Volumen/Precio.../Origen columns, full rows
# A:R) gets pushed down by 2 rows to 344:461, and two brand-new rows of data are
# inserted at the top (342:343). Using Range.Copy with an overlapping destination
# snapshots the source first (same semantics as Excel), so a single copy handles
# the whole shift - including number formats (the date column D keeps its
# date style) - in one step.
$ws.Range("A342:R459").Copy($ws.Range("A344:R461"))

# New row 342 (Primera) values
$ws.Range("D342").Value = 44524
$ws.Range("J342").Value = 2550
$ws.Range("K342").Value = 500
$ws.Range("L342").Value = 550
$ws.Range("M342").Value = 526
$ws.Range("P342").Value = 526

# New row 343 (Segunda) values
$ws.Range("D343").Value = 44524
$ws.Range("J343").Value = 1300
$ws.Range("K343").Value = 400
$ws.Range("L343").Value = 400
$ws.Range("M343").Value = 400
$ws.Range("P343").Value = 400
